# Update cryptos list data (prices and volume changes) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = '26.708.16'
$ws.Range("E2").Value = '  +0.65%  '

# Row 3
$ws.Range("D3").Value = '1.726.69'
$ws.Range("E3").Value = '  -0.37%  '

# Row 4
$ws.Range("D4").Value = '0.9974'
$ws.Range("E4").Value = '  -0.35%  '

# Row 5
$ws.Range("D5").Value = '241.76'
$ws.Range("E5").Value = '  -1.19%  '

# Row 6
$ws.Range("E6").Value = '  -0.37%  '

# Row 7
$ws.Range("D7").Value = '0.4923'
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").Value = '0.2616'
$ws.Range("E8").Value = '  -0.83%  '

# Row 9
$ws.Range("D9").Value = '0.06234'
$ws.Range("E9").Value = '  +0.75%  '

# Row 10
$ws.Range("D10").Value = '1.726.16'
$ws.Range("E10").Value = '  -0.43%  '

# Row 11
$ws.Range("D11").Value = '15.88'
$ws.Range("E11").Value = '  +1.60%  '

# Row 12
$ws.Range("D12").Value = '0.06994'
$ws.Range("E12").Value = '  -0.43%  '

# Row 13
$ws.Range("D13").Value = '0.6121'
$ws.Range("E13").Value = '  +1.57%  '

# Row 14
$ws.Range("D14").Value = '4.513'
$ws.Range("E14").Value = '  -1.12%  '

# Row 15
$ws.Range("D15").Value = '77.28'
$ws.Range("E15").Value = '  -0.34%  '

# Row 16
$ws.Range("D16").Value = '0.9979'
$ws.Range("E16").Value = '  -0.39%  '

# Row 17
$ws.Range("D17").Value = '26.684.33'
$ws.Range("E17").Value = '  +0.53%  '

# Row 18
$ws.Range("D18").Value = '0.9973'
$ws.Range("E18").Value = '  -0.38%  '

# Row 19
$ws.Range("D19").Value = '''0.000007208'
$ws.Range("E19").Value = '  +0.76%  '

# Row 20
$ws.Range("D20").Value = '11.42'
$ws.Range("E20").Value = '  -0.22%  '

# Row 21
$ws.Range("D21").Value = '1.951.44'
$ws.Range("E21").Value = '  -0.59%  '

# Row 22
$ws.Range("D22").Value = '''4.450'
$ws.Range("E22").Value = '  -1.18%  '

# Row 23
$ws.Range("D23").Value = '8.602'
$ws.Range("E23").Value = '  -0.17%  '

# Row 24
$ws.Range("D24").Value = '5.108'
$ws.Range("E24").Value = '  -1.80%  '

# Row 25
$ws.Range("D25").Value = '138.27'
$ws.Range("E25").Value = '  -0.42%  '

# Row 26
$ws.Range("D26").Value = '15.38'
$ws.Range("E26").Value = '  +0.54%  '

# Row 27
$ws.Range("D27").Value = '1.756'
$ws.Range("E27").Value = '  +1.68%  '

# Row 28
$ws.Range("D28").Value = '1.386'
$ws.Range("E28").Value = '  -3.90%  '

# Row 29
$ws.Range("D29").Value = '106.31'
$ws.Range("E29").Value = '  -0.95%  '

# Row 30
$ws.Range("D30").Value = '''3.920'
$ws.Range("E30").Value = '  -1.49%  '

# Row 31
$ws.Range("D31").Value = '0.07994'
$ws.Range("E31").Value = '  +0.15%  '

# Row 32
$ws.Range("D32").Value = '3.675'
$ws.Range("E32").Value = '  -0.16%  '

# Row 33
$ws.Range("D33").Value = '0.04502'
$ws.Range("E33").Value = '  -0.78%  '

# Row 34
$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").Value = '0.9972'
$ws.Range("E34").Value = '  -0.36%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.607'
$ws.Range("E35").Value = '  -0.35%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.003'
$ws.Range("E36").Value = '  -0.22%  '

# Row 37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.6273'
$ws.Range("E37").Value = '  +0.06%  '

# Row 38
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '0.9346'
$ws.Range("E38").Value = '  +2.89%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.022'
$ws.Range("E39").Value = '  +0.18%  '

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.415'
$ws.Range("E40").Value = '  +0.15%  '

# Row 41
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '0.9974'
$ws.Range("E41").Value = '  -0.58%  '

# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.01517'
$ws.Range("E42").Value = '  +1.63%  '

# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.599'
$ws.Range("E43").Value = '  +2.49%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '99.64'
$ws.Range("E44").Value = '  -0.97%  '

# Row 45
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.3867'
$ws.Range("E45").Value = '  -0.36%  '

# Row 46
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '6.911'
$ws.Range("E46").Value = '  +3.05%  '

# Row 47
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '''0.1160'
$ws.Range("E47").Value = '  -0.05%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05382'
$ws.Range("E48").Value = '  +0.27%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '7.866'
$ws.Range("E49").Value = '  +1.85%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '30.39'
$ws.Range("E50").Value = '  +0.01%  '

# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '51.74'
$ws.Range("E51").Value = '  +1.04%  '

Write-Output "Updated cryptos data."
